$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1908.3673
$ws.Range("I40").Value = 1724.0588
$ws.Range("J40").Value = 2006.2812
$ws.Range("K40").Value = 1724.0588
$ws.Range("L40").Value = 2006.2812
$ws.Range("M40").Value = -1549.0588
$ws.Range("N40").Value = -2356.2812
$ws.Range("H51").Value = 2742.8572
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 3033.3333
$ws.Range("K51").Value = 1000
$ws.Range("L51").Value = 3033.3333
$ws.Range("M51").Value = -516
$ws.Range("N51").Value = -4001.3333
$ws.Range("H53").Value = 213.70589
$ws.Range("I53").Value = 71.09999999999999
$ws.Range("J53").Value = 417.42856
$ws.Range("K53").Value = 71.09999999999999
$ws.Range("L53").Value = 417.42856
$ws.Range("M53").Value = 565.9
$ws.Range("N53").Value = -1691.42856
$ws.Range("H92").Value = 3315.8572
$ws.Range("I92").Value = 3355.5386
$ws.Range("K92").Value = 3355.5386
$ws.Range("M92").Value = -2107.5386
$ws.Range("H137").Value = 1126.4048
$ws.Range("I137").Value = 905.80554
$ws.Range("K137").Value = 2717.41662
$ws.Range("M137").Value = -167.41662

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3984353.5
$ws.Range("I32").Value = 4345726
$ws.Range("J32").Value = 9257.143
$ws.Range("K32").Value = 4345726
$ws.Range("L32").Value = 9257.143
$ws.Range("M32").Value = -4345439
$ws.Range("N32").Value = -9831.143
$ws.Range("H61").Value = 730.1539
$ws.Range("I61").Value = 449.2
$ws.Range("J61").Value = 1666.6666
$ws.Range("K61").Value = 449.2
$ws.Range("L61").Value = 1666.6666
$ws.Range("M61").Value = -237.2
$ws.Range("N61").Value = -2090.6666
$ws.Range("H74").Value = 1053.5927
$ws.Range("I74").Value = 1085.48
$ws.Range("K74").Value = 1085.48
$ws.Range("M74").Value = -211.48
$ws.Range("H77").Value = 1053.5927
$ws.Range("I77").Value = 1085.48
$ws.Range("K77").Value = 5427.4
$ws.Range("M77").Value = -1059.4
$ws.Range("H132").Value = 2012.4
$ws.Range("I132").Value = 1536.5555
$ws.Range("J132").Value = 3618.375
$ws.Range("K132").Value = 4609.666499999999
$ws.Range("L132").Value = 10855.125
$ws.Range("M132").Value = -2079.666499999999
$ws.Range("N132").Value = -15915.125
$ws.Range("H136").Value = 730.1539
$ws.Range("I136").Value = 449.2
$ws.Range("J136").Value = 1666.6666
$ws.Range("K136").Value = 1347.6
$ws.Range("L136").Value = 4999.9998
$ws.Range("M136").Value = 1202.4
$ws.Range("N136").Value = -10099.9998
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("H20").Value = 26335354
$ws.Range("I20").Value = 27048.13
$ws.Range("J20").Value = 66674750
$ws.Range("K20").Value = 27048.13
$ws.Range("L20").Value = 66674750
$ws.Range("M20").Value = -26801.13
$ws.Range("N20").Value = -66675244
$ws.Range("H22").Value = 317.35
$ws.Range("I22").Value = 317.35
$ws.Range("K22").Value = 317.35
$ws.Range("M22").Value = -144.35
$ws.Range("H105").Value = 3962.6667
$ws.Range("I105").Value = 4149.231
$ws.Range("K105").Value = 4149.231
$ws.Range("M105").Value = -2402.231
$ws.Range("H107").Value = 5445.1035
$ws.Range("I107").Value = 926.2222
$ws.Range("J107").Value = 66450
$ws.Range("K107").Value = 926.2222
$ws.Range("L107").Value = 66450
$ws.Range("M107").Value = 993.7778
$ws.Range("N107").Value = -70290
$ws.Range("H134").Value = 22953.229
$ws.Range("I134").Value = 1932.1143
$ws.Range("J134").Value = 79548.53999999999
$ws.Range("K134").Value = 5796.3429
$ws.Range("L134").Value = 238645.62
$ws.Range("M134").Value = -3261.3429
$ws.Range("N134").Value = -243715.62
$ws.Range("N13").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2896.2341
$ws.Range("I31").Value = 2983.75
$ws.Range("K31").Value = 2983.75
$ws.Range("M31").Value = -2688.75
$ws.Range("H34").Value = 2896.2341
$ws.Range("I34").Value = 2983.75
$ws.Range("K34").Value = 2983.75
$ws.Range("M34").Value = -2781.75
$ws.Range("H99").Value = 2081.1667
$ws.Range("I99").Value = 1882.1154
$ws.Range("K99").Value = 1882.1154
$ws.Range("M99").Value = -384.1153999999999
$ws.Range("H126").Value = 2081.1667
$ws.Range("I126").Value = 1882.1154
$ws.Range("K126").Value = 5646.3462
$ws.Range("M126").Value = -3176.3462

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 17857490
$ws.Range("I97").Value = 17857490
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 53572470
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -53571974
$ws.Range("H113").Value = 665.3333
$ws.Range("I113").Value = 796
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 2388
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -218
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 1251518.6
$ws.Range("J122").Value = 1430164.1
$ws.Range("L122").Value = 12871476.9
$ws.Range("N122").Value = -12876376.9
$ws.Range("H123").Value = 2915.7144
$ws.Range("J123").Value = 5763.3335
$ws.Range("L123").Value = 17290.0005
$ws.Range("N123").Value = -22190.0005
$ws.Range("H131").Value = 830.09
$ws.Range("I131").Value = 345
$ws.Range("J131").Value = 861.05316
$ws.Range("K131").Value = 1035
$ws.Range("L131").Value = 2583.15948
$ws.Range("M131").Value = 4005
$ws.Range("N131").Value = -12663.15948
$ws.Range("N97").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2801.1538
$ws.Range("I97").Value = 2960.7368
$ws.Range("K97").Value = 2960.7368
$ws.Range("M97").Value = -2464.7368
$ws.Range("H102").Value = 2340.8845
$ws.Range("I102").Value = 1698.0526
$ws.Range("J102").Value = 4085.7144
$ws.Range("K102").Value = 1698.0526
$ws.Range("L102").Value = 4085.7144
$ws.Range("M102").Value = -76.05259999999998
$ws.Range("N102").Value = -7329.7144
$ws.Range("H113").Value = 3591.75
$ws.Range("I113").Value = 3781.2
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 3781.2
$ws.Range("L113").Value = 750
$ws.Range("M113").Value = -1611.2
$ws.Range("H126").Value = 2041.5172
$ws.Range("I126").Value = 2033.0588
$ws.Range("J126").Value = 2053.5
$ws.Range("K126").Value = 6099.1764
$ws.Range("L126").Value = 6160.5
$ws.Range("M126").Value = -3629.1764
$ws.Range("N126").Value = -11100.5
$ws.Range("H132").Value = 3214.4443
$ws.Range("I132").Value = 2876.8333
$ws.Range("J132").Value = 3889.6667
$ws.Range("K132").Value = 8630.499899999999
$ws.Range("L132").Value = 11669.0001
$ws.Range("M132").Value = -6100.499899999999
$ws.Range("N132").Value = -16729.0001
$ws.Range("N113").Value = -5090

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9261618
$ws.Range("I7").Value = 2853.75
$ws.Range("J7").Value = 27779146
$ws.Range("K7").Value = 2853.75
$ws.Range("L7").Value = 27779146
$ws.Range("M7").Value = -2741.75
$ws.Range("N7").Value = -27779370
$ws.Range("H40").Value = 375204.2
$ws.Range("I40").Value = 460109.7
$ws.Range("K40").Value = 460109.7
$ws.Range("M40").Value = -459973.7
$ws.Range("H61").Value = 7256
$ws.Range("I61").Value = 9127.846
$ws.Range("J61").Value = 1172.5
$ws.Range("K61").Value = 9127.846
$ws.Range("L61").Value = 1172.5
$ws.Range("M61").Value = -8925.846
$ws.Range("N61").Value = -1576.5
$ws.Range("H113").Value = 7256
$ws.Range("I113").Value = 9127.846
$ws.Range("J113").Value = 1172.5
$ws.Range("K113").Value = 9127.846
$ws.Range("L113").Value = 1172.5
$ws.Range("M113").Value = -6957.846
$ws.Range("N113").Value = -5512.5
$ws.Range("H126").Value = 9261618
$ws.Range("I126").Value = 2853.75
$ws.Range("J126").Value = 27779146
$ws.Range("K126").Value = 8561.25
$ws.Range("L126").Value = 83337438
$ws.Range("M126").Value = -6091.25
$ws.Range("N126").Value = -83342378
$ws.Range("H136").Value = 2484.95
$ws.Range("I136").Value = 1608.4231
$ws.Range("J136").Value = 4112.7856
$ws.Range("K136").Value = 4825.2693
$ws.Range("L136").Value = 12338.3568
$ws.Range("M136").Value = -2275.2693
$ws.Range("N136").Value = -17438.3568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2206.5
$ws.Range("I81").Value = 1214.75
$ws.Range("J81").Value = 2999.9
$ws.Range("K81").Value = 2429.5
$ws.Range("L81").Value = 5999.8
$ws.Range("M81").Value = -1368.5
$ws.Range("N81").Value = -8121.8
$ws.Range("H84").Value = 2206.5
$ws.Range("I84").Value = 1214.75
$ws.Range("J84").Value = 2999.9
$ws.Range("K84").Value = 12147.5
$ws.Range("L84").Value = 29999
$ws.Range("M84").Value = -6843.5
$ws.Range("N84").Value = -40607
$ws.Range("H119").Value = 5000000
$ws.Range("J119").Value = 5000000
$ws.Range("L119").Value = 5000000
$ws.Range("N119").Value = -5009676
$ws.Range("H132").Value = 1702.5526
$ws.Range("I132").Value = 1213.76
$ws.Range("J132").Value = 2642.5386
$ws.Range("K132").Value = 3641.28
$ws.Range("L132").Value = 7927.6158
$ws.Range("M132").Value = -1111.28
$ws.Range("N132").Value = -12987.6158
$ws.Range("H140").Value = 47272.9
$ws.Range("J140").Value = 47272.9
$ws.Range("L140").Value = 47272.9
$ws.Range("N140").Value = -57632.9
